$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 19 de Septiembre de 2020 a las 11:53"
$ws.Range("B18").Value = 347372
$ws.Range("C18").Value = 1567
$ws.Range("D18").Value = 254386
$ws.Range("E18").Value = 88073
$ws.Range("G18").Value = 32
$ws.Range("H18").Value = 4913
$ws.Range("A47").Value = "Polonia"
$ws.Range("B47").Value = 78330
$ws.Range("C47").Value = 1002
$ws.Range("D47").Value = 63861
$ws.Range("E47").Value = 12187
$ws.Range("G47").Value = 12
$ws.Range("H47").Value = 2282
$ws.Range("A48").Value = "Japon"
$ws.Range("B48").Value = 77494
$ws.Range("D48").Value = 69899
$ws.Range("E48").Value = 6113
$ws.Range("H48").Value = 1482
$ws.Range("B98").Value = 10167
$ws.Range("C98").Value = 20
$ws.Range("D98").Value = 9315
$ws.Range("E98").Value = 722
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 130
$ws.Range("A102").Value = "Tunez"
$ws.Range("B102").Value = 9110
$ws.Range("C102").Value = 540
$ws.Range("D102").Value = 2366
$ws.Range("E102").Value = 6606
$ws.Range("G102").Value = 5
$ws.Range("H102").Value = 138
$ws.Range("A103").Value = "Finlandia"
$ws.Range("B103").Value = 8858
$ws.Range("D103").Value = 7700
$ws.Range("E103").Value = 819
$ws.Range("H103").Value = 339
$ws.Range("A104").Value = "Gabon"
$ws.Range("B104").Value = 8696
$ws.Range("D104").Value = 7848
$ws.Range("E104").Value = 795
$ws.Range("H104").Value = 53
$ws.Range("A105").Value = "Haiti"
$ws.Range("B105").Value = 8600
$ws.Range("D105").Value = 6363
$ws.Range("E105").Value = 2016
$ws.Range("H105").Value = 221
$ws.Range("A117").Value = "Hong Kong"
$ws.Range("B117").Value = 5010
$ws.Range("C117").Value = 13
$ws.Range("D117").Value = 4707
$ws.Range("E117").Value = 200
$ws.Range("H117").Value = 103
$ws.Range("A118").Value = "Cuba"
$ws.Range("B118").Value = 5004
$ws.Range("D118").Value = 4249
$ws.Range("E118").Value = 644
$ws.Range("H118").Value = 111
$ws.Range("A119").Value = "Guinea Ecuatorial"
$ws.Range("B119").Value = 5002
$ws.Range("D119").Value = 4509
$ws.Range("E119").Value = 410
$ws.Range("H119").Value = 83
$ws.Range("B128").Value = 4309
$ws.Range("C128").Value = 114
$ws.Range("D128").Value = 2981
$ws.Range("E128").Value = 1187
$ws.Range("G128").Value = 1
$ws.Range("H128").Value = 141
$ws.Range("A131").Value = "Lituania"
$ws.Range("B131").Value = 3664
$ws.Range("C131").Value = 99
$ws.Range("D131").Value = 2197
$ws.Range("E131").Value = 1380
$ws.Range("H131").Value = 87
$ws.Range("A132").Value = "Trinidad yTobago"
$ws.Range("B132").Value = 3651
$ws.Range("D132").Value = 1586
$ws.Range("E132").Value = 2005
$ws.Range("H132").Value = 60
$ws.Range("A204").Value = "Timor Oriental"
$ws.Range("A205").Value = "Santa Lucia"
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
